$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder details
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long card/account number that must stay TEXT (Excel would
# otherwise silently reinterpret a plain numeric-looking assignment as a
# number) while keeping its original cell style untouched. Stage the
# text value in a scratch cell formatted as Text, then paste only the
# value into B3 so the destination's formatting/style is left alone.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2570314725427075"
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 23.03.2024"

# Transaction row 6
$ws.Range("B6").Value = "24.03."
$ws.Range("C6").Value = "25.03."
$ws.Range("D6").Value = "PAYPAL CLARRI"
$ws.Range("E6").Value = "75,64-"

# Transaction row 7
$ws.Range("B7").Value = "27.03."
$ws.Range("C7").Value = "28.03."
$ws.Range("D7").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E7").Value = "81,01-"

# Transaction row 8
$ws.Range("B8").Value = "28.03."
$ws.Range("C8").Value = "29.03."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,55-"

# Transaction row 9 (was empty, now a new transaction)
$ws.Range("B9").Value = "31.03."
$ws.Range("C9").Value = "01.04."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-8993609"
$ws.Range("E9").Value = "55,94-"
# E9 previously used the "empty amount" style (centered); the amount
# column cells use right-aligned, non-wrapping text - line it up with
# the other amount cells (E6/E7/E8/E12).
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E9").VerticalAlignment = -4107
$ws.Range("E9").WrapText = $false

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 03.04.2024"
$ws.Range("E12").Value = "237,14-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 13.04.2024"
